$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-21 Tuesday" "2025-01-22 Wednesday"

Replace-Text "811×2=" "732×5="
Replace-Text "358×4=" "932×6="
Replace-Text "134×3=" "583×7="
Replace-Text "185×3=" "520×7="
Replace-Text "468×4=" "433×7="
Replace-Text "620×9=" "401×4="
Replace-Text "417×9=" "894×4="
Replace-Text "157×2=" "933×5="
Replace-Text "518×3=" "200×5="
Replace-Text "567×4=" "196×4="
Replace-Text "421×8=" "427×2="
Replace-Text "792×5=" "225×3="
Replace-Text "206×5=" "929×4="
Replace-Text "977×4=" "953×2="
Replace-Text "675×7=" "391×5="
Replace-Text "872×5=" "583×8="
Replace-Text "620×4=" "980×6="
Replace-Text "434×6=" "780×4="
Replace-Text "945×4=" "807×9="
Replace-Text "113×2=" "401×9="
Replace-Text "987×3=" "806×5="
Replace-Text "227×6=" "908×8="
Replace-Text "792×9=" "158×8="
Replace-Text "483×2=" "279×5="
Replace-Text "756×5=" "350×6="
